$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 355.5
$ws.Range("I6").Value = 383.8889
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 1151.6667
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = -1039.6667
$ws.Range("N6").Value = -524
$ws.Range("H28").Value = 856.4375
$ws.Range("I28").Value = 860.53845
$ws.Range("J28").Value = 838.6667
$ws.Range("K28").Value = 860.53845
$ws.Range("L28").Value = 838.6667
$ws.Range("M28").Value = -375.53845
$ws.Range("N28").Value = -1808.6667
$ws.Range("H31").Value = 94.25
$ws.Range("I31").Value = 94.25
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 282.75
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -52.75
$ws.Range("H39").Value = 48.25
$ws.Range("I39").Value = 48.25
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 144.75
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 151.25
$ws.Range("N39").Value = $null
$ws.Range("H42").Value = 12.666667
$ws.Range("I42").Value = 9
$ws.Range("J42").Value = 20
$ws.Range("K42").Value = 27
$ws.Range("L42").Value = 60
$ws.Range("M42").Value = 203
$ws.Range("N42").Value = -520
$ws.Range("H64").Value = 3875.0833
$ws.Range("I64").Value = 3954.7273
$ws.Range("J64").Value = 2999
$ws.Range("K64").Value = 3954.7273
$ws.Range("L64").Value = 2999
$ws.Range("M64").Value = -3706.7273
$ws.Range("N64").Value = -3495
$ws.Range("H67").Value = 3875.0833
$ws.Range("I67").Value = 3954.7273
$ws.Range("J67").Value = 2999
$ws.Range("K67").Value = 3954.7273
$ws.Range("L67").Value = 2999
$ws.Range("M67").Value = -3096.7273
$ws.Range("N67").Value = -4715
$ws.Range("H132").Value = 1447.5294
$ws.Range("I132").Value = 1447.5294
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4342.5882
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1812.5882

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1985.0714
$ws.Range("I88").Value = 1657.6666
$ws.Range("J88").Value = 2230.625
$ws.Range("K88").Value = 1657.6666
$ws.Range("L88").Value = 2230.625
$ws.Range("M88").Value = -1251.6666
$ws.Range("N88").Value = -3042.625
$ws.Range("H91").Value = 1985.0714
$ws.Range("I91").Value = 1657.6666
$ws.Range("J91").Value = 2230.625
$ws.Range("K91").Value = 1657.6666
$ws.Range("L91").Value = 2230.625
$ws.Range("M91").Value = -253.6666
$ws.Range("N91").Value = -5038.625
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12000
$ws.Range("N122").Value = -16900

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1598.5
$ws.Range("I31").Value = 1598.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1598.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1303.5
$ws.Range("H34").Value = 1598.5
$ws.Range("I34").Value = 1598.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1598.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1396.5
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = $null
$ws.Range("N58").Value = $null
$ws.Range("H86").Value = 126792.664
$ws.Range("I86").Value = 187139.67
$ws.Range("J86").Value = 6098.6665
$ws.Range("K86").Value = 187139.67
$ws.Range("L86").Value = 6098.6665
$ws.Range("M86").Value = -186016.67
$ws.Range("H89").Value = 126792.664
$ws.Range("I89").Value = 187139.67
$ws.Range("J89").Value = 6098.6665
$ws.Range("K89").Value = 935698.3500000001
$ws.Range("L89").Value = 30493.3325
$ws.Range("M89").Value = -930082.3500000001
$ws.Range("H99").Value = 3479.8333
$ws.Range("I99").Value = 3479.8333
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3479.8333
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1981.8333
$ws.Range("H107").Value = 1145.5
$ws.Range("I107").Value = 1145.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1145.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 774.5
$ws.Range("H126").Value = 3479.8333
$ws.Range("I126").Value = 3479.8333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10439.4999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7969.499899999999
$ws.Range("H134").Value = 1600
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 4800
$ws.Range("N134").Value = -9870
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = $null
$ws.Range("N136").Value = $null

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 300
$ws.Range("I50").Value = 300
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 900
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -419
$ws.Range("H51").Value = 450
$ws.Range("I51").Value = 450
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 1350
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -890
$ws.Range("H53").Value = 300
$ws.Range("I53").Value = 300
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 900
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -419
$ws.Range("H107").Value = 2100
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2100
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 6300
$ws.Range("N107").Value = -10140
$ws.Range("H109").Value = 10066.333
$ws.Range("I109").Value = 7599.5
$ws.Range("J109").Value = 15000
$ws.Range("K109").Value = 22798.5
$ws.Range("L109").Value = 45000
$ws.Range("M109").Value = -21758.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 1362.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1362.5
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1362.5
$ws.Range("M23").Value = $null
$ws.Range("N23").Value = -1808.5
$ws.Range("H113").Value = 5348.8887
$ws.Range("I113").Value = 5017.5
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 5017.5
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -2847.5
$ws.Range("N113").Value = -12340
$ws.Range("H122").Value = 26310.785
$ws.Range("I122").Value = 33059.363
$ws.Range("J122").Value = 1566
$ws.Range("K122").Value = 99178.08899999999
$ws.Range("L122").Value = 4698
$ws.Range("M122").Value = -96728.08899999999
$ws.Range("N122").Value = -9598
$ws.Range("H132").Value = 3905.6667
$ws.Range("I132").Value = 3608.75
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 10826.25
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -8296.25
$ws.Range("H135").Value = 80000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 80000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2313.7
$ws.Range("I61").Value = 2348.6667
$ws.Range("J61").Value = 1999
$ws.Range("K61").Value = 2348.6667
$ws.Range("L61").Value = 1999
$ws.Range("M61").Value = -2146.6667
$ws.Range("H113").Value = 2313.7
$ws.Range("I113").Value = 2348.6667
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 2348.6667
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = -178.6667000000002
$ws.Range("H122").Value = 7374.8335
$ws.Range("I122").Value = 5833.1665
$ws.Range("J122").Value = 8916.5
$ws.Range("K122").Value = 17499.4995
$ws.Range("L122").Value = 26749.5
$ws.Range("M122").Value = -15049.4995
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -12450

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 5224.75
$ws.Range("I4").Value = 10000
$ws.Range("J4").Value = 3633
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 3633
$ws.Range("M4").Value = -9887
$ws.Range("N4").Value = -3859
$ws.Range("H81").Value = 1937.25
$ws.Range("I81").Value = 874.5
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 1749
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -688
$ws.Range("H84").Value = 1937.25
$ws.Range("I84").Value = 874.5
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 8745
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -3441
$ws.Range("H107").Value = 453.81818
$ws.Range("I107").Value = 261.25
$ws.Range("J107").Value = 967.3333
$ws.Range("K107").Value = 783.75
$ws.Range("L107").Value = 2901.9999
$ws.Range("M107").Value = 1136.25
$ws.Range("N107").Value = -6741.9999
$ws.Range("H122").Value = 1726.9565
$ws.Range("I122").Value = 1765.619
$ws.Range("J122").Value = 1321
$ws.Range("K122").Value = 5296.857
$ws.Range("L122").Value = 3963
$ws.Range("M122").Value = -2846.857
$ws.Range("N122").Value = -8863
$ws.Range("H132").Value = 2908.7
$ws.Range("I132").Value = 1700
$ws.Range("J132").Value = 3043
$ws.Range("K132").Value = 5100
$ws.Range("L132").Value = 9129
$ws.Range("M132").Value = -2570
$ws.Range("N132").Value = -14189
$ws.Range("H136").Value = 1381.8572
$ws.Range("I136").Value = 749.5
$ws.Range("J136").Value = 2225
$ws.Range("K136").Value = 2248.5
$ws.Range("L136").Value = 6675
$ws.Range("M136").Value = 301.5
